# Updated cryptos list with GitHub Actions
# Applies per-cell text updates while preserving each cell's original
# "text" storage type (many values look numeric, e.g. "0.571" or
# "59.209.58", but the source workbook stores them as inline strings).
#
# Excel's COM Range.Value setter auto-coerces strings that parse as
# numbers into real numeric cells. To keep them as text (matching the
# target XML, which keeps t="inlineStr" with no style change) we:
#   1. Force the cell's number format to Text ("@") before assigning,
#      so Excel stores the literal text instead of parsing it.
#   2. Assign the value.
#   3. Clear the format again so the cell's style index reverts to the
#      sheet's default (0) instead of leaving behind a "Text format"
#      style that wasn't in the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        $Sheet,
        [string]$Addr,
        [string]$Text
    )
    $range = $Sheet.Range($Addr)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.ClearFormats()
}

Set-CellText $ws "D2" "59.209.58"
Set-CellText $ws "E2" "  -6.14%  "
Set-CellText $ws "D3" "2.460.25"
Set-CellText $ws "E3" "  -8.46%  "
Set-CellText $ws "E4" "  -0.05%  "
Set-CellText $ws "D5" "541.44"
Set-CellText $ws "E5" "  -2.73%  "
Set-CellText $ws "D6" "148.40"
Set-CellText $ws "E6" "  -6.77%  "
Set-CellText $ws "E7" "  -0.22%  "
Set-CellText $ws "D8" "0.571"
Set-CellText $ws "E8" "  -3.72%  "
Set-CellText $ws "D9" "2.478.13"
Set-CellText $ws "E9" "  -7.98%  "
Set-CellText $ws "D10" "0.0994"
Set-CellText $ws "E10" "  -6.44%  "
Set-CellText $ws "E11" "  -2.53%  "
Set-CellText $ws "D12" "5.31"
Set-CellText $ws "E12" "  -1.41%  "
Set-CellText $ws "D13" "0.353"
Set-CellText $ws "E13" "  -4.92%  "
Set-CellText $ws "D14" "2.893.26"
Set-CellText $ws "E14" "  -8.49%  "
Set-CellText $ws "D15" "24.12"
Set-CellText $ws "E15" "  -9.18%  "
Set-CellText $ws "D16" "59.116.49"
Set-CellText $ws "E16" "  -6.11%  "
Set-CellText $ws "D17" "0.0000138"
Set-CellText $ws "E17" "  -6.42%  "
Set-CellText $ws "D18" "2.527.99"
Set-CellText $ws "E18" "  -5.93%  "
Set-CellText $ws "D19" "11.16"
Set-CellText $ws "E19" "  -7.20%  "
Set-CellText $ws "D20" "4.36"
Set-CellText $ws "E20" "  -5.92%  "
Set-CellText $ws "D21" "324.90"
Set-CellText $ws "E21" "  -6.26%  "
Set-CellText $ws "D22" "0.969"
Set-CellText $ws "E22" "  -3.29%  "
Set-CellText $ws "D23" "5.76"
Set-CellText $ws "E23" "  -8.68%  "
Set-CellText $ws "D24" "0.461"
Set-CellText $ws "E24" "  -10.20%  "
Set-CellText $ws "D25" "60.74"
Set-CellText $ws "E25" "  -4.21%  "
Set-CellText $ws "E26" "  -4.68%  "
Set-CellText $ws "E27" "  -2.22%  "
Set-CellText $ws "D28" "7.73"
Set-CellText $ws "E28" "  -6.39%  "
Set-CellText $ws "B29" "Fetch.AI"
Set-CellText $ws "C29" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-CellText $ws "D29" "1.30"
Set-CellText $ws "E29" "  -8.95%  "
Set-CellText $ws "B30" "PancakeSwap"
Set-CellText $ws "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-CellText $ws "D30" "1.83"
Set-CellText $ws "E30" "  -6.57%  "
Set-CellText $ws "B31" "PEPE"
Set-CellText $ws "C31" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-CellText $ws "D31" "0.0₃0775"
Set-CellText $ws "E31" "  -10.51%  "
Set-CellText $ws "B32" "Aptos"
Set-CellText $ws "C32" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-CellText $ws "D32" "6.72"
Set-CellText $ws "E32" "  -7.76%  "
Set-CellText $ws "D33" "0.997"
Set-CellText $ws "E33" "  -0.11%  "
Set-CellText $ws "D34" "158.14"
Set-CellText $ws "E34" "  -4.44%  "
Set-CellText $ws "D35" "4.55"
Set-CellText $ws "E35" "  -8.20%  "
Set-CellText $ws "D36" "1.38"
Set-CellText $ws "E36" "  -7.62%  "
Set-CellText $ws "D37" "18.44"
Set-CellText $ws "E37" "  -5.72%  "
Set-CellText $ws "D38" "1.75"
Set-CellText $ws "E38" "  -2.54%  "
Set-CellText $ws "B39" "Bittensor"
Set-CellText $ws "C39" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-CellText $ws "D39" "320.24"
Set-CellText $ws "E39" "  -10.68%  "
Set-CellText $ws "B40" "RenderToken"
Set-CellText $ws "C40" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText $ws "D40" "5.91"
Set-CellText $ws "E40" "  -7.89%  "
Set-CellText $ws "D41" "36.58"
Set-CellText $ws "E41" "  -5.04%  "
Set-CellText $ws "E42" "  -12.98%  "
Set-CellText $ws "E43" "  -7.80%  "
Set-CellText $ws "E44" "  -0.29%  "
Set-CellText $ws "D45" "10.74"
Set-CellText $ws "E45" "  -2.74%  "
Set-CellText $ws "D46" "0.586"
Set-CellText $ws "E46" "  -5.43%  "
Set-CellText $ws "D47" "0.0943"
Set-CellText $ws "E47" "  -3.30%  "
Set-CellText $ws "D48" "0.0526"
Set-CellText $ws "E48" "  -6.86%  "
Set-CellText $ws "B49" "InjectiveProtocol"
Set-CellText $ws "C49" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-CellText $ws "D49" "19.08"
Set-CellText $ws "E49" "  -9.34%  "
Set-CellText $ws "B50" "VeChain"
Set-CellText $ws "C50" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText $ws "D50" "0.0230"
Set-CellText $ws "E50" "  -5.73%  "
Set-CellText $ws "D51" "18.54"
Set-CellText $ws "E51" "  -9.49%  "
